# Apply update: add 3 new products (rows), update totals and timestamp

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $origFmt = $range.NumberFormat
    $range.NumberFormat = "@"
    $range.Value2 = $value
    $range.NumberFormat = $origFmt
}

# Insert 3 new rows before the totals row (row 16) to make room for
# 3 new products; everything below shifts down by 3 rows.
$ws.Rows.Item(16).Insert()
$ws.Rows.Item(16).Insert()
$ws.Rows.Item(16).Insert()

# New rows 16-18 need the same look as the other data rows: copy the
# cell formatting from row 15 (last original data row) into them.
$ws.Range("A15:Q15").Copy()
$ws.Range("A16:Q18").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Recreate the merged cells for the 3 new rows (A:B, C:G, H:K, L:M, N:O)
foreach ($r in 16..18) {
    $ws.Range("A" + $r + ":B" + $r).Merge()
    $ws.Range("C" + $r + ":G" + $r).Merge()
    $ws.Range("H" + $r + ":K" + $r).Merge()
    $ws.Range("L" + $r + ":M" + $r).Merge()
    $ws.Range("N" + $r + ":O" + $r).Merge()
}

# Row heights matching the final layout
$ws.Rows.Item(7).RowHeight = 25.5
$ws.Rows.Item(8).RowHeight = 24.75
$ws.Rows.Item(9).RowHeight = 25.5
$ws.Rows.Item(10).RowHeight = 24.75
$ws.Rows.Item(11).RowHeight = 25.5
$ws.Rows.Item(12).RowHeight = 25.5
$ws.Rows.Item(13).RowHeight = 24.75
$ws.Rows.Item(14).RowHeight = 25.5
$ws.Rows.Item(15).RowHeight = 24.75
$ws.Rows.Item(16).RowHeight = 25.5
$ws.Rows.Item(17).RowHeight = 25.5
$ws.Rows.Item(18).RowHeight = 24.75

# Fill in all data rows (existing rows shift position because of the
# 2 items inserted before GLUCOVANCE and 1 before the Arabic items)
$ws.Range("A7").Value2 = 1
Set-TextValue $ws.Range("C7") "ANTODINE 20MG 30 F.C.TAB"
Set-TextValue $ws.Range("H7") "1:1"
Set-TextValue $ws.Range("L7") "1"
Set-TextValue $ws.Range("N7") "60.00"
Set-TextValue $ws.Range("P7") "19.8000"
Set-TextValue $ws.Range("Q7") "0:1"

$ws.Range("A8").Value2 = 2
Set-TextValue $ws.Range("C8") "DEPO-PEN 1.2 MIU VIAL."
Set-TextValue $ws.Range("H8") "6:0"
Set-TextValue $ws.Range("L8") "1"
Set-TextValue $ws.Range("N8") "25.00"
Set-TextValue $ws.Range("P8") "25.0000"
Set-TextValue $ws.Range("Q8") "1:0"

$ws.Range("A9").Value2 = 3
Set-TextValue $ws.Range("C9") "DEXAMETHASONE-MUP 8MG/2ML 5 AMP"
Set-TextValue $ws.Range("H9") "3:4"
Set-TextValue $ws.Range("L9") "1"
Set-TextValue $ws.Range("N9") "65.00"
Set-TextValue $ws.Range("P9") "13.0000"
Set-TextValue $ws.Range("Q9") "0:1"

$ws.Range("A10").Value2 = 4
Set-TextValue $ws.Range("C10") "DOLIPRANE 1 GM 15 TABS."
Set-TextValue $ws.Range("H10") "4:0"
Set-TextValue $ws.Range("L10") "1"
Set-TextValue $ws.Range("N10") "48.00"
Set-TextValue $ws.Range("P10") "15.8400"
Set-TextValue $ws.Range("Q10") "0:1"

$ws.Range("A11").Value2 = 5
Set-TextValue $ws.Range("C11") "GLUCOVANCE 500/5MG 30 F.C.TAB."
Set-TextValue $ws.Range("H11") "1:0"
Set-TextValue $ws.Range("L11") "1"
Set-TextValue $ws.Range("N11") "74.00"
Set-TextValue $ws.Range("P11") "74.0000"
Set-TextValue $ws.Range("Q11") "1:0"

$ws.Range("A12").Value2 = 6
Set-TextValue $ws.Range("C12") "PANADOL EXTRA 48 TAB"
Set-TextValue $ws.Range("H12") "2:1"
Set-TextValue $ws.Range("L12") "0"
Set-TextValue $ws.Range("N12") "108.00"
Set-TextValue $ws.Range("P12") "27.0000"
Set-TextValue $ws.Range("Q12") "0:1"

$ws.Range("A13").Value2 = 7
Set-TextValue $ws.Range("C13") "RIVO 320MG 20*10 TABS"
Set-TextValue $ws.Range("H13") "0:12"
Set-TextValue $ws.Range("L13") "1"
Set-TextValue $ws.Range("N13") "141.00"
Set-TextValue $ws.Range("P13") "14.1000"
Set-TextValue $ws.Range("Q13") "0:2"

$ws.Range("A14").Value2 = 8
Set-TextValue $ws.Range("C14") "SAFETRIUM 30 F.C. TABS."
Set-TextValue $ws.Range("H14") "0:3"
Set-TextValue $ws.Range("L14") "1"
Set-TextValue $ws.Range("N14") "90.00"
Set-TextValue $ws.Range("P14") "90.0000"
Set-TextValue $ws.Range("Q14") "1:0"

$ws.Range("A15").Value2 = 9
Set-TextValue $ws.Range("C15") "VOLTAREN 75MG/3ML 6 AMP."
Set-TextValue $ws.Range("H15") "0:2"
Set-TextValue $ws.Range("L15") "1"
Set-TextValue $ws.Range("N15") "102.00"
Set-TextValue $ws.Range("P15") "16.3200"
Set-TextValue $ws.Range("Q15") "0:1"

$ws.Range("A16").Value2 = 10
Set-TextValue $ws.Range("C16") "ريكسونا بليه بودر دراي"
Set-TextValue $ws.Range("H16") "2:0"
Set-TextValue $ws.Range("L16") "0"
Set-TextValue $ws.Range("N16") "37.00"
Set-TextValue $ws.Range("P16") "37.0000"
Set-TextValue $ws.Range("Q16") "1:0"

$ws.Range("A17").Value2 = 11
Set-TextValue $ws.Range("C17") "زيت فاتيكا وسط 90 مل"
Set-TextValue $ws.Range("H17") "7:0"
Set-TextValue $ws.Range("L17") "0"
Set-TextValue $ws.Range("N17") "25.00"
Set-TextValue $ws.Range("P17") "25.0000"
Set-TextValue $ws.Range("Q17") "1:0"

$ws.Range("A18").Value2 = 12
Set-TextValue $ws.Range("C18") "سرنجات 5 سم"
Set-TextValue $ws.Range("H18") "0:0"
Set-TextValue $ws.Range("L18") "0"
Set-TextValue $ws.Range("N18") "3.00"
Set-TextValue $ws.Range("P18") "6.0000"
Set-TextValue $ws.Range("Q18") "2:0"

# Update grand total (sum of sale prices) and report timestamp
$ws.Range("P19").Value2 = 363.06
Set-TextValue $ws.Range("A20") "Wednesday, 16 July, 2025 10:25 AM"

